$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column (D) number format so new rows match existing date formatting
$dateFmt = $ws.Cells.Item(2, 4).NumberFormat

function Set-Row($r, $date, $calidad, $vol, $pmin, $pmax, $pprom, $unidad, $origen, $pkg, $kgunidad) {
    $ws.Cells.Item($r, 4).NumberFormat = $dateFmt
    $ws.Cells.Item($r, 4).Value2 = $date
    $ws.Cells.Item($r, 12).Value2 = $calidad
    $ws.Cells.Item($r, 13).Value2 = $vol
    $ws.Cells.Item($r, 14).Value2 = $pmin
    $ws.Cells.Item($r, 15).Value2 = $pmax
    $ws.Cells.Item($r, 16).Value2 = $pprom
    $ws.Cells.Item($r, 17).Value2 = $unidad
    $ws.Cells.Item($r, 18).Value2 = $origen
    $ws.Cells.Item($r, 19).Value2 = $pkg
    $ws.Cells.Item($r, 20).Value2 = $kgunidad
}

Set-Row 595 44753 "Especial" 220 12500 12500 12500 "`$/caja 18 kilos" "Provincia de Curicó" 694 18
Set-Row 596 44753 "Primera" 300 11000 11000 11000 "`$/caja 18 kilos" "Provincia de Curicó" 611 18
Set-Row 597 44753 "Segunda" 280 9000 9000 9000 "`$/caja 18 kilos" "Provincia de Curicó" 500 18
Set-Row 598 44693 "Especial" 280 19800 19800 19800 "`$/caja 18 kilos" "Provincia de Curicó" 1100 18
Set-Row 599 44693 "Primera" 350 16200 16200 16200 "`$/caja 18 kilos" "Provincia de Curicó" 900 18
Set-Row 600 44693 "Segunda" 300 12600 12600 12600 "`$/caja 18 kilos" "Provincia de Curicó" 700 18
Set-Row 601 44505 "Especial" 300 11000 11000 11000 "`$/bandeja 10 kilos" "Provincia de Curicó" 1100 10
Set-Row 602 44505 "Extra (doble especial)" 280 12000 12000 12000 "`$/bandeja 10 kilos" "Provincia de Curicó" 1200 10
Set-Row 603 44505 "Primera" 220 9000 9000 9000 "`$/bandeja 10 kilos" "Provincia de Curicó" 900 10
Set-Row 604 44505 "Segunda" 350 7000 7000 7000 "`$/bandeja 10 kilos" "Provincia de Curicó" 700 10
Set-Row 605 44340 "Especial" 50 15000 15000 15000 "`$/caja 15 kilos granel" "Región Metropolitana" 1000 15
Set-Row 606 44340 "Primera" 80 13000 13000 13000 "`$/caja 15 kilos granel" "Región Metropolitana" 867 15
Set-Row 607 44340 "Segunda" 95 10000 10000 10000 "`$/caja 15 kilos granel" "Región Metropolitana" 667 15
Set-Row 608 44326 "Especial" 55 21000 21000 21000 "`$/caja 18 kilos" "Región Metropolitana" 1167 18
Set-Row 609 44326 "Primera" 80 18000 18000 18000 "`$/caja 18 kilos" "Región Metropolitana" 1000 18
Set-Row 610 44326 "Segunda" 95 14000 14000 14000 "`$/caja 18 kilos" "Región Metropolitana" 778 18
Set-Row 611 44714 "Especial" 250 10000 10000 10000 "`$/bandeja 10 kilos" "Provincia de Curicó" 1000 10
Set-Row 612 44714 "Primera" 300 8000 8000 8000 "`$/bandeja 10 kilos" "Provincia de Curicó" 800 10
Set-Row 613 44714 "Segunda" 280 6000 6000 6000 "`$/bandeja 10 kilos" "Provincia de Curicó" 600 10
Set-Row 614 44343 "Especial" 4 280000 280000 280000 "`$/bins (450 kilos)" "Región de O'Higgins" 622 450
Set-Row 615 44343 "Especial" 240 14000 15000 14583 "`$/caja 15 kilos" "Región Metropolitana" 972 15
Set-Row 616 44343 "Primera" 6 250000 250000 250000 "`$/bins (450 kilos)" "Región de O'Higgins" 556 450
Set-Row 617 44343 "Primera" 470 12000 13000 12596 "`$/caja 15 kilos" "Región Metropolitana" 840 15
Set-Row 618 44343 "Segunda" 8 200000 200000 200000 "`$/bins (450 kilos)" "Región de O'Higgins" 444 450
Set-Row 619 44343 "Segunda" 410 9000 9500 9293 "`$/caja 15 kilos" "Región Metropolitana" 620 15
Set-Row 620 44426 "Especial" 75 9000 9000 9000 "`$/bandeja 10 kilos" "Región de O'Higgins" 900 10
Set-Row 621 44426 "Extra (doble especial)" 50 11000 11000 11000 "`$/bandeja 10 kilos" "Región de O'Higgins" 1100 10
Set-Row 622 44426 "Primera" 130 7000 7000 7000 "`$/bandeja 10 kilos" "Región de O'Higgins" 700 10
Set-Row 623 44376 "Especial" 100 14400 14400 14400 "`$/caja 18 kilos" "Región de O'Higgins" 800 18
Set-Row 624 44376 "Extra (doble especial)" 60 19800 19800 19800 "`$/caja 18 kilos" "Región de O'Higgins" 1100 18
Set-Row 625 44376 "Primera" 190 10800 10800 10800 "`$/caja 18 kilos" "Región de O'Higgins" 600 18
Set-Row 626 44406 "Especial" 180 8000 8000 8000 "`$/bandeja 10 kilos" "Región Metropolitana" 800 10
Set-Row 627 44406 "Extra (doble especial)" 110 10000 10000 10000 "`$/bandeja 10 kilos" "Región Metropolitana" 1000 10
Set-Row 628 44406 "Primera" 250 6000 6000 6000 "`$/bandeja 10 kilos" "Región Metropolitana" 600 10
Set-Row 629 44454 "Especial" 110 10000 10000 10000 "`$/bandeja 10 kilos" "Provincia de Curicó" 1000 10
Set-Row 630 44454 "Extra (doble especial)" 70 12000 12000 12000 "`$/bandeja 10 kilos" "Provincia de Curicó" 1200 10
Set-Row 631 44454 "Primera" 150 8000 8000 8000 "`$/bandeja 10 kilos" "Provincia de Curicó" 800 10

# Fill in the constant leading columns (A-C, E-K) for the 3 newly appended rows
foreach ($r in 629..631) {
    $ws.Cells.Item($r, 1).Value2 = 9
    $ws.Cells.Item($r, 2).Value2 = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($r, 3).Value2 = "Metropolitana"
    $ws.Cells.Item($r, 5).Value2 = 13
    $ws.Cells.Item($r, 6).Value2 = "Fruta"
    $ws.Cells.Item($r, 7).Value2 = 100101
    $ws.Cells.Item($r, 8).Value2 = "Berries"
    $ws.Cells.Item($r, 9).Value2 = 100101007
    $ws.Cells.Item($r, 10).Value2 = "Kiwi"
    $ws.Cells.Item($r, 11).Value2 = "Hayward"
}
